$d = $word.ActiveDocument

# Locate the paragraph that ends with "Adivina adivinador". Searching via
# Find leaves $searchRange collapsed around the matched text, which we use
# as an anchor instead of hard-coding character offsets.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Adivina adivinador", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Position right after the last character of "Adivina adivinador" but
    # still inside the existing paragraph (i.e. before its paragraph mark,
    # where the _GoBack bookmark currently sits).
    $insertPoint = $d.Range($searchRange.End, $searchRange.End)

    # Insert the new sentence's text first. Because the bookmark sits at the
    # very end of the paragraph, typing here places the new text before the
    # bookmark, exactly like pasting new content right before it.
    $insertPoint.InsertAfter("No conformidad en la pegada")

    # Now split the paragraph right before the text we just typed. This turns
    # "No conformidad en la pegada" (together with the trailing bookmark)
    # into its own new paragraph, leaving "Adivina adivinador" alone in the
    # original paragraph.
    $breakPoint = $d.Range($searchRange.End, $searchRange.End)
    $breakPoint.InsertParagraphAfter()
}
